$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-6 from 45204 to 45207
$ws.Range("C2:C6").Value = 45207
